$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: Validation -> F1 train
$ws.Range("O1").Value = "F1 train"

# Update O column (Validation / F1 train metric) values for rows 2-16
$ws.Range("O2").Value = 0.8461538461538461
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("O5").Value = 0.7654320987654321
$ws.Range("O6").Value = 0.59375
$ws.Range("O7").Value = 0.925
$ws.Range("O8").Value = 0.9736842105263158
$ws.Range("O9").Value = 0.8493150684931506
$ws.Range("O10").Value = 1
$ws.Range("O11").Value = 0.6440677966101694
$ws.Range("O12").Value = 0.935064935064935
$ws.Range("O13").Value = 0.9333333333333333
$ws.Range("O14").Value = 1
$ws.Range("O15").Value = 0.868421052631579

# Row 16 (MLP, Free) updates
$ws.Range("C16").Value = "{'activation': 'relu', 'alpha': 0.0001, 'hidden_layer_sizes': (64,), 'learning_rate': 'constant'}"
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 6
$ws.Range("I16").Value = 0.55
$ws.Range("J16").Value = 0.4705882352941176
$ws.Range("K16").Value = 0.4
$ws.Range("L16").Value = 0.5714285714285714
$ws.Range("M16").Value = 0.7
$ws.Range("N16").Value = 0.4
$ws.Range("O16").Value = 0.7466666666666667
